$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @{
  2  = @(3.457262170408171, 1.007702210158641, 0.3492156199738474, 1.343180230842464, 6.463794665469322, 0.0007918837209675011, 0.009995295743722019, 0.002670933310996659)
  3  = @(2.997853015548628, 0.869853693225707, 0.3064966316520383, 1.155751310872887, 5.654258087503649, 0.0007995267261627653, 0.00566132105441941, 0.0006838291904367466)
  4  = @(2.717878933285135, 0.7870470786523356, 0.2802437474734631, 1.042418103720038, 5.159655295851422, 0.00080432863316101, 0.003562661497822361, 0.0004342247631967666)
  5  = @(2.604150420224698, 0.7552167163809997, 0.2683125171909211, 0.9965441341777392, 4.944718871251297, 0.0008063365970838446, 0.00283109922356406, 0.000622890166188661)
  6  = @(2.585207879043878, 0.7518079961090791, 0.2648516623235935, 0.9888719871438241, 4.892255946361189, 0.0008066988388415279, 0.002711539273524188, 0.0007586302123687716)
  7  = @(2.716100113317282, 0.7916813407236987, 0.2760110113271139, 1.041580167724106, 5.110605349627804, 0.0008044293528601984, 0.00353482695453966, 0.0006621355425160402)
  8  = @(3.297930760598206, 0.9664754913100921, 0.3289524093292329, 1.277802295241244, 6.121670233424283, 0.0007945947571169861, 0.00834023203055928, 0.001998973943032389)
  9  = @(4.461276162754871, 1.317628976875312, 0.4409129175473083, 1.758791490145597, 8.221729476999883, 0.0007759650060087116, 0.02293966575955908, 0.01318788530029646)
  10 = @(5.332556541617578, 1.589096750276838, 0.5003877811189454, 2.021280822810638, 9.52940342855976, 0.0007632755018147897, 0.03654587955837574, 0.02786243332301552)
  11 = @(5.696287717332439, 1.68929629314573, 0.3475534396946358, 1.319623735495057, 7.984121445527933, 0.0007625219121048961, 0.05102710592101189, 0.02986244529849014)
  12 = @(5.821453138132654, 1.71081905899149, 0.2378295168972926, 0.8095450099342685, 6.621540193137264, 0.0007639042150792652, 0.08537059510579326, 0.02834500312373667)
  13 = @(5.769062860697261, 1.681431827745541, 0.1489218976295632, 0.4192333535295205, 5.257492853800528, 0.0007670939994285801, 0.1363447470369579, 0.02441139653021018)
  14 = @(5.653258874752623, 1.640501665964337, 0.09914364752189897, 0.2200208786317432, 4.320288365004302, 0.0007700778533965794, 0.18241339872894, 0.02096641771609153)
  15 = @(5.587027398541977, 1.621696695691526, 0.08790478637317278, 0.180148105521539, 4.065077586257132, 0.0007712711082461959, 0.1938160907104987, 0.01972158660848322)
  16 = @(5.224085017559617, 1.516492032032318, 0.08526566423042681, 0.1700411051112454, 3.828971836814645, 0.00077591099244641, 0.1772398316277872, 0.01491254195543057)
  17 = @(5.007352928785792, 1.457340512143503, 0.1063005318635462, 0.2527426306022065, 4.119620197415884, 0.0007779609381707517, 0.1380845818779335, 0.01298819220147251)
  18 = @(4.891181072940014, 1.426784906721878, 0.1630917184647416, 0.488672800946965, 4.989782952390868, 0.0007777175351106891, 0.08624958060428156, 0.0128610273449894)
  19 = @(4.865410637195282, 1.433989826580387, 0.2585924094822332, 0.9316599478342198, 6.311405617291229, 0.0007753471943131185, 0.04486607482094485, 0.01491509726212481)
  20 = @(5.100185226391659, 1.532472938516833, 0.4702190910753643, 1.945769640609583, 9.0203272447026, 0.0007668416119520997, 0.03244554102579755, 0.02379408360536761)
  21 = @(5.784892976581773, 1.750257957264921, 0.5534414014092022, 2.318371252717682, 10.45960704905036, 0.000756201365715177, 0.0464444912688986, 0.03888522512563775)
  22 = @(6.238143054835405, 1.888477988335353, 0.6034973891750326, 2.513470465521181, 11.35561038052936, 0.0007494717359220136, 0.05622256578920748, 0.05032041316143676)
  23 = @(5.996341369730999, 1.808243448853432, 0.5820632899846885, 2.409407682906746, 10.93547678253219, 0.0007529736297253577, 0.05097258498324031, 0.04393539907997379)
  24 = @(5.093396454990398, 1.522929633599688, 0.494095960497134, 2.024357664267967, 9.282577626360364, 0.0007664096614360714, 0.03328234609949288, 0.02388743652415215)
  25 = @(4.142656613778001, 1.229542174013773, 0.4026513398619613, 1.625711697670212, 7.560761790821061, 0.0007810655889403182, 0.01826779157226266, 0.009293847328576632)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  $ws.Cells.Item($r, 2).Value = $vals[0]
  $ws.Cells.Item($r, 3).Value = $vals[1]
  $ws.Cells.Item($r, 4).Value = $vals[2]
  $ws.Cells.Item($r, 5).Value = $vals[3]
  $ws.Cells.Item($r, 6).Value = $vals[4]
  $ws.Cells.Item($r, 7).Value = $vals[5]
  $ws.Cells.Item($r, 8).Value = $vals[6]
  $ws.Cells.Item($r, 9).Value = $vals[7]
  $ws.Cells.Item($r, 16).Value = 0
  $ws.Cells.Item($r, 17).Value = 0
}
